# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Ixion_Profits leve-profit tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2974.875
$ws.Range("J62").Value = 2287.375
$ws.Range("L62").Value = 2287.375
$ws.Range("N62").Value = -3535.375
$ws.Range("H65").Value = 2974.875
$ws.Range("J65").Value = 2287.375
$ws.Range("L65").Value = 11436.875
$ws.Range("N65").Value = -17676.875
$ws.Range("H92").Value = 129630300
$ws.Range("I92").Value = 55555556
$ws.Range("J92").Value = 166667680
$ws.Range("K92").Value = 55555556
$ws.Range("L92").Value = 166667680
$ws.Range("M92").Value = -55554308
$ws.Range("N92").Value = -166670176
$ws.Range("H112").Value = 47619988
$ws.Range("J112").Value = 63493184
$ws.Range("L112").Value = 190479552
$ws.Range("N112").Value = -190481768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1448.1765
$ws.Range("I2").Value = 1172.8667
$ws.Range("J2").Value = 3513
$ws.Range("K2").Value = 1172.8667
$ws.Range("L2").Value = 3513
$ws.Range("M2").Value = -1059.8667
$ws.Range("N2").Value = -3739
$ws.Range("H4").Value = 261.25
$ws.Range("I4").Value = 178
$ws.Range("K4").Value = 178
$ws.Range("M4").Value = -62
$ws.Range("H32").Value = 2015.59
$ws.Range("I32").Value = 2015.59
$ws.Range("K32").Value = 2015.59
$ws.Range("M32").Value = -1728.59
$ws.Range("H61").Value = 462542.12
$ws.Range("I61").Value = 13257.3
$ws.Range("K61").Value = 13257.3
$ws.Range("M61").Value = -13045.3
$ws.Range("H62").Value = 42300
$ws.Range("J62").Value = 42300
$ws.Range("L62").Value = 42300
$ws.Range("N62").Value = -43548
$ws.Range("H65").Value = 42300
$ws.Range("J65").Value = 42300
$ws.Range("L65").Value = 126900
$ws.Range("N65").Value = -133140
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H110").Value = 1470.3334
$ws.Range("I110").Value = 1470.3334
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1470.3334
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 574.6666
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 1448.1765
$ws.Range("I116").Value = 1172.8667
$ws.Range("J116").Value = 3513
$ws.Range("K116").Value = 1172.8667
$ws.Range("L116").Value = 3513
$ws.Range("M116").Value = 1121.1333
$ws.Range("N116").Value = -8101
$ws.Range("H132").Value = 2567271.8
$ws.Range("I132").Value = 2177.5938
$ws.Range("K132").Value = 6532.7814
$ws.Range("M132").Value = -4002.7814
$ws.Range("H136").Value = 462542.12
$ws.Range("I136").Value = 13257.3
$ws.Range("K136").Value = 39771.89999999999
$ws.Range("M136").Value = -37221.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1448.1765
$ws.Range("I3").Value = 1172.8667
$ws.Range("J3").Value = 3513
$ws.Range("K3").Value = 1172.8667
$ws.Range("L3").Value = 3513
$ws.Range("M3").Value = -1058.8667
$ws.Range("N3").Value = -3741
$ws.Range("H7").Value = 251.5
$ws.Range("I7").Value = 251.5
$ws.Range("K7").Value = 251.5
$ws.Range("M7").Value = -138.5
$ws.Range("H94").Value = 1371.6786
$ws.Range("I94").Value = 822.4737
$ws.Range("J94").Value = 2531.111
$ws.Range("K94").Value = 822.4737
$ws.Range("L94").Value = 2531.111
$ws.Range("M94").Value = -371.4737
$ws.Range("N94").Value = -3433.111
$ws.Range("H99").Value = 40001404
$ws.Range("I99").Value = 52632870
$ws.Range("J99").Value = 1749.8334
$ws.Range("K99").Value = 52632870
$ws.Range("L99").Value = 1749.8334
$ws.Range("M99").Value = -52631372
$ws.Range("N99").Value = -4745.8334
$ws.Range("H105").Value = 5423.6924
$ws.Range("I105").Value = 6026
$ws.Range("J105").Value = 4460
$ws.Range("K105").Value = 6026
$ws.Range("L105").Value = 4460
$ws.Range("M105").Value = -4279
$ws.Range("N105").Value = -7954
$ws.Range("H134").Value = 20802.828
$ws.Range("I134").Value = 3868.628
$ws.Range("K134").Value = 11605.884
$ws.Range("M134").Value = -9070.884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2052.375
$ws.Range("I12").Value = 1641.6
$ws.Range("J12").Value = 2737
$ws.Range("K12").Value = 1641.6
$ws.Range("L12").Value = 2737
$ws.Range("M12").Value = -1471.6
$ws.Range("N12").Value = -3077
$ws.Range("H94").Value = 7686.0713
$ws.Range("I94").Value = 7181.2
$ws.Range("J94").Value = 7966.5557
$ws.Range("K94").Value = 7181.2
$ws.Range("L94").Value = 7966.5557
$ws.Range("M94").Value = -6730.2
$ws.Range("N94").Value = -8868.555700000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2858007.2
$ws.Range("I131").Value = 5556135
$ws.Range("J131").Value = 1165.9412
$ws.Range("K131").Value = 16668405
$ws.Range("L131").Value = 3497.8236
$ws.Range("M131").Value = -16663365
$ws.Range("N131").Value = -13577.8236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1442.2122
$ws.Range("I102").Value = 1341.7097
$ws.Range("K102").Value = 1341.7097
$ws.Range("M102").Value = 280.2902999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 43690.5
$ws.Range("J81").Value = 43690.5
$ws.Range("L81").Value = 43690.5
$ws.Range("N81").Value = -45686.5
$ws.Range("H84").Value = 43690.5
$ws.Range("J84").Value = 43690.5
$ws.Range("L84").Value = 131071.5
$ws.Range("N84").Value = -141055.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 30173
$ws.Range("J76").Value = 30173
$ws.Range("L76").Value = 30173
$ws.Range("N76").Value = -30803
$ws.Range("H79").Value = 30173
$ws.Range("J79").Value = 30173
$ws.Range("L79").Value = 30173
$ws.Range("N79").Value = -32357
$ws.Range("H122").Value = 1920.6154
$ws.Range("I122").Value = 1406.6316
$ws.Range("J122").Value = 3315.7144
$ws.Range("K122").Value = 4219.8948
$ws.Range("L122").Value = 9947.143199999999
$ws.Range("M122").Value = -1769.8948
$ws.Range("N122").Value = -14847.1432
